{"js": "// box_label_template: append \". {{stage}}\" right after \" {{object}}\" in the\n// \"\u041e\u0431\u044a\u0435\u043a\u0442:\" table cell, as its own run with the same Times New Roman / 36\n// half-pt / en-US character formatting as the \" {{object}}\" run.\n\nconst body = context.document.body;\n\n// Locate the \" {{object}}\" placeholder text.\nconst results = body.search(\"{{object}}\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find '{{object}}' placeholder in the document\");\n}\n\nconst objectRange = results.items[0];\n\n// Insert the new text immediately after the found range (same paragraph,\n// same cell). insertText() inherits the formatting of the text it is\n// anchored to, so the new text starts out merged into the neighbouring run.\nconst insertedRange = objectRange.insertText(\". {{stage}}\", Word.InsertLocation.after);\n\n// Force the engine to materialize the inserted text as its own run (instead\n// of silently merging it back into the \" {{object}}\" run because the\n// formatting is identical) by nudging the font size away and back. This\n// keeps the original rFonts (ascii/hAnsi/cs = \"Times New Roman\") intact on\n// both runs.\ninsertedRange.font.size = 20;\nawait context.sync();\n\ninsertedRange.font.size = 18;\nawait context.sync();\n", "ps1": "# box_label_template: append \". {{stage}}\" right after \" {{object}}\" in the\n# \"\u041e\u0431\u044a\u0435\u043a\u0442:\" table cell, as its own run with the same Times New Roman / 36\n# half-pt / en-US character formatting as the \" {{object}}\" run.\n\n$d = $word.ActiveDocument\n\n# Locate the \" {{object}}\" run and collapse the range to its end so the new\n# text is inserted immediately after it (same paragraph, same cell).\n$rng = $d.Content\n$found = $rng.Find.Execute(\"{{object}}\")\nif (-not $found) {\n    throw \"Could not find '{{object}}' placeholder in the document\"\n}\n$rng.Collapse(0)\n$rng.InsertAfter(\". {{stage}}\")\n\n# InsertAfter() re-adopts the inherited formatting of the preceding run, so\n# the inserted text is merged into the neighbouring run. Re-select just the\n# newly inserted text (fresh Find so the range handle isn't stale after the\n# structural edit) and nudge its size away-and-back; that forces the engine\n# to materialize the run as a distinct <w:r> while keeping the original\n# rFonts (ascii/hAnsi/cs=\"Times New Roman\") intact.\n$newRng = $d.Content\n$newRng.Find.Execute(\". {{stage}}\")\n$newRng.Font.Size = 20\n$newRng.Font.Size = 18\n"}
